# updated bionwire sequence gen system to be deterministic
#
# For each data sheet (Sequences, Names, Descriptions), row 8 (plate position "G")
# holds 7 generated BioNWire binding-handle values in columns B..H. The
# deterministic regeneration reshuffles which value lands in which column
# according to a fixed permutation (same on every sheet):
#   new B <- old C      new E <- old G
#   new C <- old F      new F <- old B
#   new D <- old H      new G <- old D
#   new H <- old E
$wb = $excel.ActiveWorkbook

function Permute-RowEight {
    param($ws)

    $b = $ws.Range("B8").Value2
    $c = $ws.Range("C8").Value2
    $d = $ws.Range("D8").Value2
    $e = $ws.Range("E8").Value2
    $f = $ws.Range("F8").Value2
    $g = $ws.Range("G8").Value2
    $h = $ws.Range("H8").Value2

    $ws.Range("B8").Value = $c
    $ws.Range("C8").Value = $f
    $ws.Range("D8").Value = $h
    $ws.Range("E8").Value = $g
    $ws.Range("F8").Value = $b
    $ws.Range("G8").Value = $d
    $ws.Range("H8").Value = $e
}

$wsSequences = $wb.Worksheets.Item("Sequences")
$wsNames = $wb.Worksheets.Item("Names")
$wsDescriptions = $wb.Worksheets.Item("Descriptions")

Permute-RowEight $wsSequences
Permute-RowEight $wsNames
Permute-RowEight $wsDescriptions

# View/selection state also shifted around as part of the edit:
# - "Descriptions" was the active tab, now "Sequences" is.
# - each sheet's lingering selection rectangle moved to a new cell.
[void]$wsNames.Range("B12").Select()

$wsDescriptions.Activate()
[void]$wsDescriptions.Range("I8").Select()

$wsSequences.Activate()
[void]$wsSequences.Range("C47").Select()
